$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old fragment caption row ("mation"/"pompes)"/"Hiver"/"Eté"/"Année")
# so the unit-label row (m3/s / MW / MW / GWh / GWh / GWh) becomes row 1.
$ws.Rows(1).Delete()

# Build the new header row 1: idx / idx2 / Name / Date Start / Date End,
# then relabel the existing unit headers to be unambiguous per-column names.
$ws.Cells.Item(1,1).Value = "idx"
$ws.Cells.Item(1,2).Value = "idx2"
$ws.Cells.Item(1,3).Value = "Name"
$ws.Cells.Item(1,4).Value = "Date Start"
$ws.Cells.Item(1,5).Value = "Date End"
$ws.Cells.Item(1,6).Value = "(m3/s)"
$ws.Cells.Item(1,7).Value = "(MW1)"
$ws.Cells.Item(1,8).Value = "(MW2)"
$ws.Cells.Item(1,9).Value = "(GWh) Winter"
$ws.Cells.Item(1,10).Value = "(GWh) Summer"
$ws.Cells.Item(1,11).Value = "(GWh) Year"

# Leave the selection on the first data row, matching the refreshed table.
[void]$ws.Range("A2:K2").Select()
